# Generate Report for Handoff
#
# A new handoff run replaces the old source-document identifier
# (0e30c74f-814a-4e95-8fda-835800d48082) with a new one
# (45042dd7-9921-4597-a615-dab45f678036), the old handoff-archive hash
# (d4dbf63cf18c05f62f1eaac67cb37d4f66546612) with a new one
# (99f410cd76057e55a9156ff3b3e73b3e493d21b5), and bumps the "Latest Handoff
# Datetime" timestamps to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$oldGuid = "0e30c74f-814a-4e95-8fda-835800d48082"
$newGuid = "45042dd7-9921-4597-a615-dab45f678036"
$oldHash = "d4dbf63cf18c05f62f1eaac67cb37d4f66546612"
$newHash = "99f410cd76057e55a9156ff3b3e73b3e493d21b5"

$newFileName = "$newGuid.md"
$newZhFile = "$newGuid.$newHash.zh-cn.xlf"
$newDeFile = "$newGuid.$newHash.de-de.xlf"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Update the visible cell text (shared strings) ------------------------
# A2 (File Name / Source File Name) is shared across all three sheets.
$wsOverview.Range("A2").Value = $newFileName
$wsZhCn.Range("A2").Value = $newFileName
$wsDeDe.Range("A2").Value = $newFileName

# zh-cn sheet: Latest Handoff File (C2) and Latest Handoff Datetime (D2)
$wsZhCn.Range("C2").Value = $newZhFile
$wsZhCn.Range("D2").Value = "2016-03-08 08:38:37"

# de-de sheet: Latest Handoff File (C2) and Latest Handoff Datetime (D2)
$wsDeDe.Range("C2").Value = $newDeFile
$wsDeDe.Range("D2").Value = "2016-03-08 08:38:40"

# --- Update the hyperlink display text to match the new names -------------
# The underlying link targets (URLs) are NOT part of this change (the diff
# touches only sharedStrings.xml text and the <hyperlink display="..."/>
# attributes) - only the text shown in the cell (the hyperlink's
# "display"/TextToDisplay) needs to track the new file names. Writing to an
# existing hyperlink's properties in place isn't supported by this host (it
# always appends a fresh entry instead of updating), so each sheet's
# hyperlinks are rebuilt from scratch: delete all, then re-add every link
# (changed and unchanged alike) with its ORIGINAL, untouched target URL and
# the up-to-date display text.

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/07da35fb3dfb2f126e8edf66b44c604358027d11/e2e/$oldGuid.md"
$configTarget = "https://github.com/OpenLocalizationTest/oltest/blob/07da35fb3dfb2f126e8edf66b44c604358027d11/.localization-config"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/017a450572b46f17d7e79106871b2101e9dbe41f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c821daba289d2c0a28078db6c2791fac482bc8d0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    $mdTarget,
    "",
    "",
    $newFileName) | Out-Null
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    $configTarget,
    "",
    "",
    ".localization-config") | Out-Null

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    $mdTarget,
    "",
    "",
    $newFileName) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C2"),
    $zhXlfTarget,
    "",
    "",
    $newZhFile) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    $configTarget,
    "",
    "",
    ".localization-config") | Out-Null

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    $mdTarget,
    "",
    "",
    $newFileName) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C2"),
    $deXlfTarget,
    "",
    "",
    $newDeFile) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    $configTarget,
    "",
    "",
    ".localization-config") | Out-Null
